$d = $word.ActiveDocument

$find = "Datas da campaña de Constelación de Pegaso 2022: 8-17 de outubro, 7-16 de novembro"
$replace = "Datas da campaña de 2022 que usan Constelación de Pegaso: 8-17 de outubro, 7-16 de novembro"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
